$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 24 - existing rows 24..112 shift down to 25..113
$ws.Rows("24:24").Insert()

# Populate the newly inserted row 24 with its data
$ws.Cells.Item(24, 1).Value = 5
$ws.Cells.Item(24, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(24, 3).Value = "Maule"
$ws.Cells.Item(24, 4).Value = 44608
$ws.Cells.Item(24, 5).Value = 7
$ws.Cells.Item(24, 6).Value = 100112030
$ws.Cells.Item(24, 7).Value = "Poroto granado"
$ws.Cells.Item(24, 8).Value = "Sin especificar"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 500
$ws.Cells.Item(24, 11).Value = 17000
$ws.Cells.Item(24, 12).Value = 17000
$ws.Cells.Item(24, 13).Value = 17000
$ws.Cells.Item(24, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(24, 15).Value = "Región del Maule"
$ws.Cells.Item(24, 16).Value = 680
$ws.Cells.Item(24, 17).Value = 25
$ws.Cells.Item(24, 18).Value = "Hortaliza"
